$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.876.04"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.04%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.444.82"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.09%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "560.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "163.75"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.509"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.37%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.171"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +8.66%  "
$ws.Range("E10").Value = "  -2.04%  "
$ws.Range("E11").Value = "  -0.32%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.60"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.86%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000180"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "68.742.18"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.10%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.890.64"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.52%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "23.41"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.05%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.441.97"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.73%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.61"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.11%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "339.29"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.13%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.04"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.51%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "3.85"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.28%  "
$ws.Range("E22").Value = "  +3.09%  "
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("E24").Value = "  -1.22%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.76"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.567.38"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.18%  "
$ws.Range("E27").Value = "  +2.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.29"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.29%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0825"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.18"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.79%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.19"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.87%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "431.96"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.82%  "
$ws.Range("E34").Value = "  -1.61%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "160.17"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.63%  "
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.04"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.80%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.106"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.34%  "
$ws.Range("E40").Value = "  +0.95%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.52"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.68%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.39"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.94%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.08"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.27%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.07"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.82%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.34"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.29%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "130.71"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.05%  "
$ws.Range("E47").Value = "  -0.31%  "
$ws.Range("E48").Value = "  -0.49%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.560"
$ws.Range("D49").Style = "Normal"
$ws.Range("E50").Value = "  +3.19%  "
$ws.Range("E51").Value = "  +0.54%  "
